$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Locate the paragraph that contains "Dr. Pepper on the local machine..."
# (the admin/local-machine walkthrough paragraph).
# -----------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $cand = $d.Paragraphs.Item($i)
  if ($cand.Range.Text -like "*Dr. Pepper on the local machine has a quantity of 12*") {
    $target = $cand
    break
  }
}

# -----------------------------------------------------------------------
# Step 1: Re-write the tail of that paragraph (from "Dr. Pepper on the
# local machine..." through to the end of the paragraph) as a clean run
# structure: one merged run for the long sentence run-on, then the
# gramStart/"has"/gramEnd proof-error wrapped run, then the trailing
# run.  This also drops the stray mid-sentence _GoBack bookmark and the
# lastRenderedPageBreak marker that used to sit in the middle of it.
# -----------------------------------------------------------------------
$findRange = $target.Range.Duplicate
$found = $findRange.Find.Execute("Dr. Pepper on the local machine has a quantity of 12", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$startPos = $findRange.Start
$paraEnd = $target.Range.End - 1   # stop just before the paragraph mark
$tailRange = $d.Range($startPos, $paraEnd)

$tailXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Dr. Pepper on the local machine has a quantity of 12.  I can also click the next button and view all items on the local machine.  The third button shows the products in order of quantity from least to greatest.  The fourth button will show me all products with a quantity of 3 or lower.  This would make it easy for purchasing addition products before they run out.  The remote machine options </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>has</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> all the same buttons and you can click on them.  These are products not on the local machine.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tailRange.InsertXML($tailXml)

# -----------------------------------------------------------------------
# Step 2: Insert six new paragraphs right after that paragraph: three
# bold+underlined headings ("Show Inventory Purchase order", "Show
# Difficult coding examples", "Run Process Queue") each followed by a
# blue tabbed note paragraph.  The _GoBack bookmark is re-created at the
# very end of the last note paragraph.
# -----------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $cand = $d.Paragraphs.Item($i)
  if ($cand.Range.Text -like "*Dr. Pepper on the local machine has a quantity of 12*") {
    $target = $cand
    break
  }
}
$insertPos = $target.Range.End - 1   # just before the paragraph mark
$insertRange = $d.Range($insertPos, $insertPos)

$newParasXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Show Inventory Purchase order</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:tab/><w:t>(Show inventory purchase order and how it works)</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Show Difficult coding examples</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:tab/><w:t>(Show difficult coding examples and explain how we fixed them, probably 2 would be good)</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Run Process Queue</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:tab/><w:t>(Run Process Queue and explain how it works and the end result) (examples Richard purchase &#8220;X&#8221; for &#8220;y&#8221;, Chad bought &#8220;X&#8221; for &#8220;y&#8221; etc.)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertRange.InsertXML($newParasXml)

Write-Host "Done. Paragraph count now:" $d.Paragraphs.Count
